{"js": "// Remove the trailing \"Ver no Jupiter ...\" line, the site footer/copyright\n// line (\"\u00a9 2020 . Contact: ...\") and the blank paragraph that separates them\n// from the \"LOQ4044: ...\" requirement line above (this mirrors the\n// Jekyll site rebuild that dropped the scraped page-chrome text).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text.trim();\n  if (targetTexts.indexOf(text) !== -1) {\n    // Also remove the immediately preceding blank paragraph (the spacer\n    // paragraph that sits right before the \"Ver no Jupiter ...\" line).\n    if (text === targetTexts[0] && i > 0 && items[i - 1].text.trim() === \"\") {\n      items[i - 1].delete();\n    }\n    items[i].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" line, the site footer/copyright\n# line (\"(c) 2020 . Contact: ...\") and the blank spacer paragraph that sits\n# right before them (right after the \"LOQ4044: ...\" requirement line).\n#\n# We locate the two text paragraphs by their content, then build a single\n# Range spanning from the blank spacer paragraph just before the first of\n# them through to the end of the copyright paragraph, and delete that whole\n# range in one shot. Deleting as a single contiguous range (rather than\n# paragraph-by-paragraph) avoids index/reference invalidation issues that\n# happen when the Paragraphs collection is mutated while it is being walked.\n\n$d = $word.ActiveDocument\n\n$footerText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = [char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$footerIndex = -1\n$copyrightIndex = -1\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq $footerText) {\n        $footerIndex = $i\n    } elseif ($text -eq $copyrightText) {\n        $copyrightIndex = $i\n    }\n}\n\nif ($footerIndex -gt 0 -and $copyrightIndex -eq ($footerIndex + 1)) {\n    $startIndex = $footerIndex\n    # Include the blank spacer paragraph right before the footer line, if any.\n    if ($footerIndex -gt 1) {\n        $prevText = $d.Paragraphs.Item($footerIndex - 1).Range.Text.TrimEnd([char]13, [char]7)\n        if ($prevText -eq \"\") {\n            $startIndex = $footerIndex - 1\n        }\n    }\n\n    $startRange = $d.Paragraphs.Item($startIndex).Range\n    $endRange = $d.Paragraphs.Item($copyrightIndex).Range\n    $fullRange = $d.Range($startRange.Start, $endRange.End)\n    $fullRange.Delete()\n}\n"}
